$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5480.2
$ws.Range("I18").Value = 645
$ws.Range("J18").Value = 8703.666999999999
$ws.Range("K18").Value = 645
$ws.Range("L18").Value = 8703.666999999999
$ws.Range("M18").Value = -361
$ws.Range("N18").Value = -9271.666999999999
$ws.Range("H64").Value = 23791.584
$ws.Range("I64").Value = 28389.111
$ws.Range("K64").Value = 28389.111
$ws.Range("M64").Value = -28141.111
$ws.Range("H67").Value = 23791.584
$ws.Range("I67").Value = 28389.111
$ws.Range("K67").Value = 28389.111
$ws.Range("M67").Value = -27531.111
$ws.Range("H127").Value = 2003
$ws.Range("I127").Value = 1170.2858
$ws.Range("K127").Value = 3510.8574
$ws.Range("M127").Value = 1449.1426
$ws.Range("H137").Value = 36782.47
$ws.Range("I137").Value = 28045.455
$ws.Range("K137").Value = 84136.36500000001
$ws.Range("M137").Value = -81586.36500000001
$ws.Range("H138").Value = 29216.648
$ws.Range("I138").Value = 1646.12
$ws.Range("J138").Value = 86655.25
$ws.Range("K138").Value = 4938.36
$ws.Range("L138").Value = 259965.75
$ws.Range("M138").Value = 201.6400000000003
$ws.Range("N138").Value = -270245.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H34").Value = 380007
$ws.Range("J34").Value = 340009.34
$ws.Range("L34").Value = 340009.34
$ws.Range("N34").Value = -340551.34
$ws.Range("H61").Value = 7580.5293
$ws.Range("I61").Value = 1373.7693
$ws.Range("K61").Value = 1373.7693
$ws.Range("M61").Value = -1161.7693
$ws.Range("H74").Value = 612413.7
$ws.Range("I74").Value = 2000670
$ws.Range("J74").Value = 17446.715
$ws.Range("K74").Value = 2000670
$ws.Range("L74").Value = 17446.715
$ws.Range("M74").Value = -1999796
$ws.Range("N74").Value = -19194.715
$ws.Range("H77").Value = 612413.7
$ws.Range("I77").Value = 2000670
$ws.Range("J77").Value = 17446.715
$ws.Range("K77").Value = 10003350
$ws.Range("L77").Value = 87233.575
$ws.Range("M77").Value = -9998982
$ws.Range("N77").Value = -95969.575
$ws.Range("H102").Value = 2473
$ws.Range("I102").Value = 2356.4285
$ws.Range("K102").Value = 2356.4285
$ws.Range("M102").Value = -734.4285
$ws.Range("H132").Value = 2194.1365
$ws.Range("I132").Value = 1913.6
$ws.Range("K132").Value = 5740.799999999999
$ws.Range("M132").Value = -3210.799999999999
$ws.Range("H136").Value = 7580.5293
$ws.Range("I136").Value = 1373.7693
$ws.Range("K136").Value = 4121.3079
$ws.Range("M136").Value = -1571.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2517.4
$ws.Range("I86").Value = 2396.75
$ws.Range("K86").Value = 2396.75
$ws.Range("M86").Value = -1273.75
$ws.Range("H89").Value = 2517.4
$ws.Range("I89").Value = 2396.75
$ws.Range("K89").Value = 11983.75
$ws.Range("M89").Value = -6367.75
$ws.Range("H102").Value = 18909.785
$ws.Range("H123").Value = 39062.5
$ws.Range("J123").Value = 39062.5
$ws.Range("L123").Value = 39062.5
$ws.Range("N123").Value = -48862.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 763.5909
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 807.61536
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 807.61536
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1507.61536
$ws.Range("H31").Value = 5883341
$ws.Range("I31").Value = 7143763
$ws.Range("J31").Value = 1369.6666
$ws.Range("K31").Value = 7143763
$ws.Range("L31").Value = 1369.6666
$ws.Range("M31").Value = -7143468
$ws.Range("N31").Value = -1959.6666
$ws.Range("H34").Value = 5883341
$ws.Range("I34").Value = 7143763
$ws.Range("J34").Value = 1369.6666
$ws.Range("K34").Value = 7143763
$ws.Range("L34").Value = 1369.6666
$ws.Range("M34").Value = -7143561
$ws.Range("N34").Value = -1773.6666
$ws.Range("H62").Value = 4460.4
$ws.Range("I62").Value = 4448
$ws.Range("K62").Value = 4448
$ws.Range("M62").Value = -3824
$ws.Range("H65").Value = 4460.4
$ws.Range("I65").Value = 4448
$ws.Range("K65").Value = 22240
$ws.Range("M65").Value = -19120
$ws.Range("H99").Value = 4964.35
$ws.Range("I99").Value = 4154.7144
$ws.Range("J99").Value = 6853.5
$ws.Range("K99").Value = 4154.7144
$ws.Range("L99").Value = 6853.5
$ws.Range("M99").Value = -2656.7144
$ws.Range("N99").Value = -9849.5
$ws.Range("H115").Value = 29999.6
$ws.Range("J115").Value = 29642.785
$ws.Range("L115").Value = 29642.785
$ws.Range("N115").Value = -31992.785
$ws.Range("H126").Value = 4964.35
$ws.Range("I126").Value = 4154.7144
$ws.Range("J126").Value = 6853.5
$ws.Range("K126").Value = 12464.1432
$ws.Range("L126").Value = 20560.5
$ws.Range("M126").Value = -9994.143199999999
$ws.Range("N126").Value = -25500.5
$ws.Range("H132").Value = 65205
$ws.Range("I132").Value = 101694
$ws.Range("K132").Value = 305082
$ws.Range("M132").Value = -302552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 499.83334
$ws.Range("I5").Value = 499.83334
$ws.Range("K5").Value = 1499.50002
$ws.Range("M5").Value = -1387.50002
$ws.Range("H10").Value = 630.8
$ws.Range("I10").Value = 388
$ws.Range("J10").Value = 995
$ws.Range("K10").Value = 1164
$ws.Range("L10").Value = 2985
$ws.Range("M10").Value = -1025
$ws.Range("N10").Value = -3263
$ws.Range("H69").Value = 8416.666999999999
$ws.Range("I69").Value = 25500
$ws.Range("K69").Value = 76500
$ws.Range("M69").Value = -75689
$ws.Range("H70").Value = 4894.737
$ws.Range("H72").Value = 8416.666999999999
$ws.Range("I72").Value = 25500
$ws.Range("K72").Value = 229500
$ws.Range("M72").Value = -225444
$ws.Range("H73").Value = 4894.737
$ws.Range("H122").Value = 985.6667
$ws.Range("I122").Value = 650
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400
$ws.Range("H135").Value = 499.83334
$ws.Range("I135").Value = 499.83334
$ws.Range("K135").Value = 4498.50006
$ws.Range("M135").Value = -1963.50006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 28333.334
$ws.Range("J74").Value = 28333.334
$ws.Range("L74").Value = 28333.334
$ws.Range("N74").Value = -30205.334
$ws.Range("H77").Value = 28333.334
$ws.Range("J77").Value = 28333.334
$ws.Range("L77").Value = 85000.00199999999
$ws.Range("N77").Value = -94360.00199999999
$ws.Range("H97").Value = 2303.1667
$ws.Range("I97").Value = 2826.6667
$ws.Range("J97").Value = 1779.6666
$ws.Range("K97").Value = 2826.6667
$ws.Range("L97").Value = 1779.6666
$ws.Range("M97").Value = -2330.6667
$ws.Range("N97").Value = -2771.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10386
$ws.Range("H55").Value = 2413.9333
$ws.Range("I55").Value = 1640.3334
$ws.Range("K55").Value = 1640.3334
$ws.Range("M55").Value = -1467.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17782.584
$ws.Range("J45").Value = 17782.584
$ws.Range("L45").Value = 17782.584
$ws.Range("N45").Value = -18764.584
$ws.Range("H121").Value = 104999.5
$ws.Range("J121").Value = 104999.5
$ws.Range("L121").Value = 104999.5
$ws.Range("N121").Value = -108493.5
$ws.Range("H132").Value = 103571.43
$ws.Range("J132").Value = 4500
$ws.Range("N132").Value = -18560
